$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

# ------------------------------------------------------------------
# Helper: write a value into a cell as TEXT, never letting Excel's
# autodetect turn a date-looking string (e.g. "2012-04-05") into a
# date serial number. Uses a scratch cell formatted as Text, then
# copies only the *value* back via a formats-paste trick so the
# destination keeps ordinary (unstyled) formatting.
# ------------------------------------------------------------------
function Set-TextValue($range, $value) {
    $scratch = $ws.Range("ZZ1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $range.PasteSpecial(-4122)   # xlPasteFormats (applies Text format)
    $range.Value = $value
    $scratch.Clear()
}

# --- Row 1: the sheet originally had a stray duplicate data row here;
#     turn it into the real column header row. ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "dealer"
$ws.Range("E1").Value = "quantity"
$ws.Range("F1").Value = "face_value"
$ws.Range("G1").Value = "currency"
$ws.Range("H1").Value = "total"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# Give the new I1:O1 header cells the same (bold/bordered) look as the
# rest of row 1.
$ws.Range("B1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)

# --- Row 2 held the wrong "total" (text "458621" instead of a number);
#     fix it, then append the record's extra fields. ---
$ws.Range("H2").Value = 458621

# --- Rows 2-6: append the per-record metadata columns. ---
$records = @(
    @{ Row = 2; Index = 82 },
    @{ Row = 3; Index = 83 },
    @{ Row = 4; Index = 84 },
    @{ Row = 5; Index = 85 },
    @{ Row = 6; Index = 86 }
)

foreach ($rec in $records) {
    $r = $rec.Row
    $ws.Range("I$r").Value = "fund"
    $ws.Range("J$r").Value = "normal"
    Set-TextValue $ws.Range("K$r") "2012-04-05"
    $ws.Range("L$r").Value = "陳鎮湘"
    $ws.Range("M$r").Value = 1754
    $ws.Range("N$r").Value = "tmpc08e1"
    $ws.Range("O$r").Value = $rec.Index
}

# Match the new I2:O6 data cells to the rest of the (unstyled) data rows.
$ws.Range("B2").Copy()
$ws.Range("I2:O6").PasteSpecial(-4122)
